# Update dSF (column F) values for several rows to reflect repulled data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    11 = 0
    13 = -1
    16 = 3
    18 = -2
    19 = -5
    25 = -1
    29 = -2
    34 = -4
    36 = -6
    38 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
